# Refresh the cryptos list (Price + Volume(1h) columns) with the latest
# scraped values. Row 47/48 also swap which coin (PaxDollar / EnergySwap)
# occupies which rank.
#
# Note: several "Price" values (column D) look like plain decimal numbers
# (e.g. "1.000", "0.9996"); Excel would normally auto-convert those to
# numeric cells. They are forced back to literal text via a leading
# apostrophe (exactly what typing '1.000 into a cell does in real Excel),
# then the cell style is reset to "Normal" so no visible quote-prefix
# formatting lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.118.55"
$ws.Range("E2").Value = "  -2.66%  "
$ws.Range("D3").Value = "1.864.25"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'306.73"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.5161"
$ws.Range("E7").Value = "  +3.20%  "
$ws.Range("D8").Value = "'0.3751"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").Value = "'0.07156"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "'0.8848"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "1.879.90"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").Value = "'0.07561"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "'5.333"
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").Value = "'89.24"
$ws.Range("E15").Value = "  -2.62%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "'0.000008551"
$ws.Range("D18").Value = "'14.15"
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "27.168.61"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("D21").Value = "'5.027"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").Value = "2.117.27"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").Value = "'10.60"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("D25").Value = "'150.97"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("D26").Value = "'1.848"
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("D27").Value = "'18.01"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").Value = "'2.156"
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("D29").Value = "'112.85"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").Value = "'4.739"
$ws.Range("E30").Value = "  -3.47%  "
$ws.Range("D31").Value = "'4.692"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").Value = "'0.09016"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "'0.05152"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").Value = "'3.103"
$ws.Range("E34").Value = "  -3.41%  "
$ws.Range("D35").Value = "'0.7545"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  -4.93%  "
$ws.Range("D37").Value = "'0.02032"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("D38").Value = "'2.534"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").Value = "'3.021"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "'1.082"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D41").Value = "'0.5353"
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("D42").Value = "'6.657"
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("D43").Value = "'114.90"
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("D44").Value = "'8.517"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "'0.1483"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").Value = "'0.4667"
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.15"
$ws.Range("E47").Value = "  -4.21%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'0.9992"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "'1.571"
$ws.Range("E49").Value = "  -3.75%  "
$ws.Range("D50").Value = "'64.89"
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("D51").Value = "'36.34"
$ws.Range("E51").Value = "  -1.81%  "

foreach ($addr in @("D4","D5","D6","D7","D8","D9","D10","D13","D14","D15","D16","D17","D18","D21","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")) {
    $ws.Range($addr).Style = "Normal"
}
